$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Enter the text via a formula that evaluates to the literal string, then
    # convert the formula result to a static value. This avoids Excel's
    # automatic conversion of date-looking text (e.g. "2025-01-14") into a
    # real date value/format while keeping the cell a plain shared string.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy($range) | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Header row
Set-TextValue $ws.Range("A1") "Inmueble"
Set-TextValue $ws.Range("B1") "fecha_apertura"
Set-TextValue $ws.Range("C1") "fecha_cierre"
Set-TextValue $ws.Range("D1") "comentario"

# Row 2
Set-TextValue $ws.Range("A2") "PRUEBA"
Set-TextValue $ws.Range("B2") "2025-01-14"
Set-TextValue $ws.Range("C2") "0000-00-00"
Set-TextValue $ws.Range("D2") "WWW"

# Row 3
Set-TextValue $ws.Range("A3") "PRUEBA"
Set-TextValue $ws.Range("B3") "2025-01-15"
Set-TextValue $ws.Range("C3") "0000-00-00"
Set-TextValue $ws.Range("D3") "asww"

# Row 4
Set-TextValue $ws.Range("A4") "PRUEBA"
Set-TextValue $ws.Range("B4") "2025-01-14"
Set-TextValue $ws.Range("C4") "2025-01-31"
Set-TextValue $ws.Range("D4") "qqqqqqqqqqq"
